$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) with refreshed data
$ws.Range("B9").Value = 4253494.15
$ws.Range("C9").Value = 672306.24
$ws.Range("D9").Value = 4925800.390000001
$ws.Range("E9").Value = 13.64866999817668
$ws.Range("F9").Value = 86.35133000182331
$ws.Range("G9").Value = -35.02495417745801
$ws.Range("H9").Value = -23.18777499911431
$ws.Range("I9").Value = 42821
$ws.Range("J9").Value = 1846
$ws.Range("K9").Value = 44667
$ws.Range("L9").Value = 31057
$ws.Range("M9").Value = 158.6051579354091
$ws.Range("N9").Value = 8.282971749443613

$wb.Save()
